# Update ticket "want to go" counts (column F) and mark one listing as
# sold out (column G) on the "展览" (Exhibitions) and "全部类型" (All
# types) sheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet "Exhibitions") ----
$wsExpo.Range("F2").Value  = 183
$wsExpo.Range("F3").Value  = 105
$wsExpo.Range("F6").Value  = 5169
$wsExpo.Range("F7").Value  = 428
$wsExpo.Range("F8").Value  = 601
$wsExpo.Range("F9").Value  = 896
$wsExpo.Range("F10").Value = 809
$wsExpo.Range("F13").Value = 552
$wsExpo.Range("F17").Value = 1713
$wsExpo.Range("F18").Value = 1439
$wsExpo.Range("F19").Value = 797
$wsExpo.Range("G20").Value = "已售罄"
$wsExpo.Range("F27").Value = 525
$wsExpo.Range("F28").Value = 2388
$wsExpo.Range("F29").Value = 167
$wsExpo.Range("F30").Value = 90
$wsExpo.Range("F31").Value = 79
$wsExpo.Range("F33").Value = 235
$wsExpo.Range("F37").Value = 226
$wsExpo.Range("F39").Value = 617
$wsExpo.Range("F43").Value = 56

# ---- 全部类型 (sheet "All types") ----
$wsAll.Range("F3").Value  = 183
$wsAll.Range("F4").Value  = 105
$wsAll.Range("F7").Value  = 5169
$wsAll.Range("F8").Value  = 428
$wsAll.Range("F9").Value  = 601
$wsAll.Range("F12").Value = 896
$wsAll.Range("F13").Value = 809
$wsAll.Range("F17").Value = 552
$wsAll.Range("F22").Value = 1713
$wsAll.Range("F23").Value = 1439
$wsAll.Range("F24").Value = 797
$wsAll.Range("G25").Value = "已售罄"
$wsAll.Range("F32").Value = 525
$wsAll.Range("F33").Value = 2388
$wsAll.Range("F34").Value = 167
$wsAll.Range("F35").Value = 90
$wsAll.Range("F36").Value = 79
$wsAll.Range("F38").Value = 235
$wsAll.Range("F43").Value = 617
$wsAll.Range("F46").Value = 56
